$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "localisation"
$ws.Range("B2").Value = "ser_pub_loc___canton"

# Delete row 3 entirely (which also shrinks the table range)
$ws.Rows("3:3").Delete()

# Restore selection to match the target state (bottom-right pane active cell B5)
$ws.Range("B5").Select() | Out-Null
